$d = $word.ActiveDocument

# Locate the paragraph that ends the HULA-025 user story block:
# "Para: Tenerlo disponible para asignarlo a reparaciones o mantenimientos de vehículos"
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Tenerlo disponible para asignarlo a reparaciones*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found"
}

$anchor = $d.Paragraphs.Item($anchorIndex).Range

# Insert the new HULA-026 user story block right after the anchor paragraph,
# before the existing blank separator paragraph.
$anchor.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 1).Range.Text = "HULA-026 Actualización de datos del Taller"

$d.Paragraphs.Item($anchorIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 2).Range.Text = "Como: Empleado del área"

$d.Paragraphs.Item($anchorIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 3).Range.Text = "Quiero: Editar o actualizar la información de un taller ya registrado (como herramientas, repuestos o ubicación)"

$d.Paragraphs.Item($anchorIndex + 3).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 4).Range.Text = "Para: Mantener la información actualizada para una mejor gestión del taller y sus recursos"

$d.Paragraphs.Item($anchorIndex + 4).Range.InsertParagraphAfter()
